$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Auto-update draw results: append the new Pick 3 row for 2025-11-30.
$row = 75

# Force the new cells to be stored as plain text (matching the rest of the
# sheet, where dates/phase codes/results are all text, not real
# dates/numbers) by switching to a text number format before writing.
$target = $ws.Range("A" + $row + ":E" + $row)
$target.NumberFormat = "@"

$ws.Range("A" + $row).Value = "2025-11-30"
$ws.Range("B" + $row).Value = "Pick 3"
$ws.Range("C" + $row).Value = "251130"
$ws.Range("D" + $row).Value = "0-4-2"
$ws.Range("E" + $row).Value = "2025-11-30T21:37:36.792+04:00"

# Restore the default "Normal" style so no stray formatting/style is left
# behind on the new row (values stay text since they're already committed).
$target.Style = "Normal"
